$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.248.52"
$ws.Range("E2").Value = '  +0.80%  '

$ws.Range("D3").Value = "'1.896.55"
$ws.Range("E3").Value = '  +0.54%  '

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = "'323.05"
$ws.Range("E5").Value = '  -2.07%  '

$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = '  -0.03%  '

$ws.Range("D7").Value = "'0.4698"
$ws.Range("E7").Value = '  +2.55%  '

$ws.Range("D8").Value = "'0.4021"
$ws.Range("E8").Value = '  -1.80%  '

$ws.Range("D9").Value = "'47.49"
$ws.Range("E9").Value = '  -0.32%  '

$ws.Range("D10").Value = "'0.08000"
$ws.Range("E10").Value = '  +0.60%  '

$ws.Range("D11").Value = "'0.9936"
$ws.Range("E11").Value = '  -0.30%  '

$ws.Range("D12").Value = "'22.48"
$ws.Range("E12").Value = '  +3.44%  '

$ws.Range("D13").Value = "'1.869.56"
$ws.Range("E13").Value = '  -2.27%  '

$ws.Range("D14").Value = "'5.853"
$ws.Range("E14").Value = '  -1.01%  '

$ws.Range("D15").Value = "'7.039"

$ws.Range("D16").Value = "'88.99"
$ws.Range("E16").Value = '  +0.56%  '

$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = '  +0.02%  '

$ws.Range("D18").Value = "'0.06622"
$ws.Range("E18").Value = '  +1.07%  '

$ws.Range("D19").Value = "'0.00001026"
$ws.Range("E19").Value = '  +0.16%  '

$ws.Range("D20").Value = "'17.48"
$ws.Range("E20").Value = '  +0.42%  '

$ws.Range("D21").Value = "'0.9997"
$ws.Range("E21").Value = '  -0.07%  '

$ws.Range("D22").Value = "'29.255.63"
$ws.Range("E22").Value = '  +0.76%  '

$ws.Range("D23").Value = "'5.499"
$ws.Range("E23").Value = '  +1.27%  '

$ws.Range("D25").Value = "'2.196"
$ws.Range("E25").Value = '  -0.21%  '

$ws.Range("D26").Value = "'2.118.80"
$ws.Range("E26").Value = '  -0.61%  '

$ws.Range("D27").Value = "'154.59"
$ws.Range("E27").Value = '  -0.83%  '

$ws.Range("D28").Value = "'19.66"
$ws.Range("E28").Value = '  +0.53%  '

$ws.Range("D29").Value = "'6.079"
$ws.Range("E29").Value = '  +10.35%  '

$ws.Range("D30").Value = "'2.084"
$ws.Range("E30").Value = '  +0.08%  '

$ws.Range("D31").Value = "'117.21"
$ws.Range("E31").Value = '  -0.13%  '

$ws.Range("D32").Value = "'1.056"
$ws.Range("E32").Value = '  +1.95%  '

$ws.Range("D33").Value = "'0.09437"
$ws.Range("E33").Value = '  +1.43%  '

$ws.Range("D34").Value = "'1.396"
$ws.Range("E34").Value = '  -0.80%  '

$ws.Range("D35").Value = "'3.538"
$ws.Range("E35").Value = '  +0.43%  '

$ws.Range("D36").Value = "'5.346"
$ws.Range("E36").Value = '  +1.15%  '

$ws.Range("D37").Value = "'0.06073"
$ws.Range("E37").Value = '  +0.38%  '

$ws.Range("D38").Value = "'0.02241"
$ws.Range("E38").Value = '  +0.80%  '

$ws.Range("D39").Value = "'1.175"
$ws.Range("E39").Value = '  +0.32%  '

$ws.Range("D40").Value = "'8.053"
$ws.Range("E40").Value = '  -3.37%  '

$ws.Range("E41").Value = '  +0.53%  '

$ws.Range("D42").Value = "'0.1825"
$ws.Range("E42").Value = '  +0.24%  '

$ws.Range("D43").Value = "'2.480"
$ws.Range("E43").Value = '  +9.09%  '

$ws.Range("D44").Value = "'10.03"
$ws.Range("E44").Value = '  -0.50%  '

$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D45").Value = "'0.07696"
$ws.Range("E45").Value = '  +2.22%  '

$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Value = "'1.248"
$ws.Range("E46").Value = '  -0.79%  '

$ws.Range("D47").Value = "'12.18"
$ws.Range("E47").Value = '  +1.63%  '

$ws.Range("D48").Value = "'0.5472"
$ws.Range("E48").Value = '  +0.37%  '

$ws.Range("D49").Value = "'1.900"
$ws.Range("E49").Value = '  +0.02%  '

$ws.Range("D50").Value = "'113.32"
$ws.Range("E50").Value = '  +1.85%  '

$ws.Range("B51").Value = 'WOONetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D51").Value = "'0.2949"
$ws.Range("E51").Value = '  +6.12%  '

$ws.Range("B2:E51").Style = "Normal"
